$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$oldText = $cellA1.Value2
$newText = $oldText.Replace(
    "1000 Bs = 1.86 = 6733.56 pesos",
    "1000 Bs = 1.85 = 6740.87 pesos"
).Replace(
    "6733.56 pesos = 1.85 = 887.03 Bs",
    "6740.87 pesos = 1.85 = 920.51 Bs"
)
$cellA1.Value2 = $newText

# --- Sheet "tasas": update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 539.99
$wsTasas.Range("O10").Value = 3640
$wsTasas.Range("N12").Value = 3639.5
$wsTasas.Range("O12").Value = 497
